$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Project Name and Neighborhood for row 3
$ws.Range("A3").Value = "Savannah Condopark"
$ws.Range("B3").Value = "Simei"

# Clear the Officer (M3) cell for row 3
$ws.Range("M3").ClearContents()

# Update visibility status from Hidden to Visible for both rows
$ws.Range("N2").Value = "Visible"
$ws.Range("N3").Value = "Visible"
